$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Trening" header in column F, matching the style of the
# existing header row (bold, bordered, centered).
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "Trening"

# Column A holds date-time serials. Build the timestamp number format on
# A2 first (through a lowercase pattern, then the uppercase one that is
# actually used), then fan that formatting out to the rest of the date
# column before filling in all the row values below.
$ws.Range("A2").Value = 45675.47986111111
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A2").Copy($ws.Range("A3:A13"))

# Rewrite the data rows (2-7 replaced, 8-13 new) with the expanded data
# set that also records which training segment ("Trening") each sample
# belongs to.
$ws.Range("A2").Value = 45675.47986111111
$ws.Range("B2").Value = 1317.7
$ws.Range("C2").Value = 11.11
$ws.Range("D2").Value = 2.80811756
$ws.Range("E2").Value = "10-15"
$ws.Range("F2").Value = "Duża Gra"

$ws.Range("A3").Value = 45675.48263888889
$ws.Range("B3").Value = 1586.6
$ws.Range("C3").Value = 14.68
$ws.Range("D3").Value = 2.828296933857143
$ws.Range("E3").Value = "10-15"
$ws.Range("F3").Value = "Duża Gra"

$ws.Range("A4").Value = 45675.48888888889
$ws.Range("B4").Value = 2081.7
$ws.Range("C4").Value = 11.41
$ws.Range("D4").Value = 2.837149858428571
$ws.Range("E4").Value = "10-15"
$ws.Range("F4").Value = "Duża Gra"

$ws.Range("A5").Value = 45675.47986111111
$ws.Range("B5").Value = 1317.5
$ws.Range("C5").Value = 8.960000000000001
$ws.Range("D5").Value = 2.567775538857143
$ws.Range("E5").Value = "5-10"
$ws.Range("F5").Value = "Duża Gra"

$ws.Range("A6").Value = 45675.48541666667
$ws.Range("B6").Value = 1784
$ws.Range("C6").Value = 9.789999999999999
$ws.Range("D6").Value = 2.558486223142857
$ws.Range("E6").Value = "5-10"
$ws.Range("F6").Value = "Duża Gra"

$ws.Range("A7").Value = 45675.48888888889
$ws.Range("B7").Value = 2081.5
$ws.Range("C7").Value = 9.609999999999999
$ws.Range("D7").Value = 2.725956082285714
$ws.Range("E7").Value = "5-10"
$ws.Range("F7").Value = "Duża Gra"

$ws.Range("A8").Value = 45675.47777777778
$ws.Range("B8").Value = 1121.6
$ws.Range("C8").Value = 10.09
$ws.Range("D8").Value = 2.542425734571428
$ws.Range("E8").Value = "10-15"
$ws.Range("F8").Value = "Mała Gra"

$ws.Range("A9").Value = 45675.47777777778
$ws.Range("B9").Value = 1125.4
$ws.Range("C9").Value = 10.11
$ws.Range("D9").Value = 2.645646912714286
$ws.Range("E9").Value = "10-15"
$ws.Range("F9").Value = "Mała Gra"

$ws.Range("A10").Value = 45675.47847222222
$ws.Range("B10").Value = 1178.1
$ws.Range("C10").Value = 12.75
$ws.Range("D10").Value = 2.931171042428571
$ws.Range("E10").Value = "10-15"
$ws.Range("F10").Value = "Mała Gra"

$ws.Range("A11").Value = 45675.47777777778
$ws.Range("B11").Value = 1121.4
$ws.Range("C11").Value = 8.44
$ws.Range("D11").Value = 2.632095745571429
$ws.Range("E11").Value = "5-10"
$ws.Range("F11").Value = "Mała Gra"

$ws.Range("A12").Value = 45675.47777777778
$ws.Range("B12").Value = 1125.3
$ws.Range("C12").Value = 9.09
$ws.Range("D12").Value = 2.612198880714285
$ws.Range("E12").Value = "5-10"
$ws.Range("F12").Value = "Mała Gra"

$ws.Range("A13").Value = 45675.47847222222
$ws.Range("B13").Value = 1177.8
$ws.Range("C13").Value = 9.539999999999999
$ws.Range("D13").Value = 2.220314068285714
$ws.Range("E13").Value = "5-10"
$ws.Range("F13").Value = "Mała Gra"
